$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[54.07260261047598, 72.99033498856672]"
$ws.Range("Q2").Value = "[1.3773949772495797, 1.7044476659115801]"
$ws.Range("U2").Value = "[43.71295320919167, 56.518566479671506]"
$ws.Range("M3").Value = "[54.0970185382852, 72.9659190607575]"
$ws.Range("U3").Value = "[43.71793877751156, 56.51358091135161]"
$ws.Range("M4").Value = "[53.1264742915145, 73.77045252205566]"
$ws.Range("N4").Value = [double]"4.440892098500626e-16"
$ws.Range("O4").Value = [double]"4.440892098500626e-16"
$ws.Range("Q4").Value = "[1.2264475824825025, 1.603816069400195]"
$ws.Range("U4").Value = "[43.76716993934688, 56.69130106008799]"
$ws.Range("Y4").Value = [double]"18.11963963963991"
$ws.Range("Z4").Value = [double]"19.58090090090119"
$ws.Range("M5").Value = "[52.735639005258776, 73.99433168112277]"
$ws.Range("N5").Value = [double]"1.332267629550188e-15"
$ws.Range("O5").Value = [double]"1.332267629550188e-15"
$ws.Range("Q5").Value = "[1.1635528346628865, 1.540921321580579]"
$ws.Range("U5").Value = "[41.71810076661528, 54.68290138164157]"
$ws.Range("Y5").Value = [double]"18.36318318318346"
$ws.Range("Z5").Value = [double]"19.82444444444474"
$ws.Range("M6").Value = "[52.769807070257926, 74.0777156292174]"
$ws.Range("N6").Value = [double]"1.332267629550188e-15"
$ws.Range("O6").Value = [double]"1.332267629550188e-15"
$ws.Range("U6").Value = "[43.83435867659661, 56.84313116451386]"
$ws.Range("M7").Value = "[51.590543466221504, 75.30989717155849]"
$ws.Range("N7").Value = [double]"4.75175454539567e-14"
$ws.Range("O7").Value = [double]"4.75175454539567e-14"
$ws.Range("U7").Value = "[43.90606553509494, 56.97846283037521]"
$ws.Range("M8").Value = "[51.24050994805495, 75.562187886311]"
$ws.Range("N8").Value = [double]"1.101341240428155e-13"
$ws.Range("O8").Value = [double]"1.101341240428155e-13"
$ws.Range("Q8").Value = "[0.8742369946926551, 1.2767633807381937]"
$ws.Range("R8").Value = [double]"4.951594689828198e-14"
$ws.Range("S8").Value = [double]"4.951594689828198e-14"
$ws.Range("U8").Value = "[41.86464435222517, 54.949066588705335]"
$ws.Range("Y8").Value = [double]"19.38606606606636"
$ws.Range("Z8").Value = [double]"20.94474474474506"
$ws.Range("M9").Value = "[51.07865433327875, 75.95702914239867]"
$ws.Range("N9").Value = [double]"2.142730437526552e-13"
$ws.Range("O9").Value = [double]"2.142730437526552e-13"
$ws.Range("U9").Value = "[43.99159984162748, 57.08449677983162]"
$ws.Range("M10").Value = "[50.877503258949226, 76.35134228367096]"
$ws.Range("N10").Value = [double]"4.318767565791859e-13"
$ws.Range("O10").Value = [double]"4.318767565791859e-13"
$ws.Range("U10").Value = "[44.06816420708641, 57.18081686911536]"
$ws.Range("M11").Value = "[51.454337348768306, 77.69954881425885]"
$ws.Range("N11").Value = [double]"6.865619184281968e-13"
$ws.Range("O11").Value = [double]"6.865619184281968e-13"
$ws.Range("U11").Value = "[44.3612517836482, 57.71377786445461]"
$ws.Range("M12").Value = "[52.90749222766715, 77.5019284336187]"
$ws.Range("N12").Value = [double]"6.372680161348399e-14"
$ws.Range("O12").Value = [double]"6.372680161348399e-14"
$ws.Range("U12").Value = "[44.68785817319232, 57.96129015335098]"
$ws.Range("M13").Value = "[54.80777263153317, 75.60167725888033]"
$ws.Range("U13").Value = "[44.69552051053814, 57.9535359114162]"
$ws.Range("M14").Value = "[53.009343114741256, 76.14445324108709]"
$ws.Range("N14").Value = [double]"1.176836406102666e-14"
$ws.Range("O14").Value = [double]"1.176836406102666e-14"
$ws.Range("U14").Value = "[44.36131017538303, 57.71415499142952]"
$ws.Range("M15").Value = "[50.404534659901145, 76.11635733819058]"
$ws.Range("N15").Value = [double]"6.87894186057747e-13"
$ws.Range("O15").Value = [double]"6.87894186057747e-13"
$ws.Range("U15").Value = "[43.535873838607294, 56.90525755602666]"
$ws.Range("M16").Value = "[50.11979872866675, 75.2631308985201]"
$ws.Range("N16").Value = [double]"4.534150832569139e-13"
$ws.Range("O16").Value = [double]"4.534150832569139e-13"
$ws.Range("Q16").Value = "[-1.0692107129334643, -0.6666843268879239]"
$ws.Range("R16").Value = [double]"3.538391801782836e-11"
$ws.Range("S16").Value = [double]"3.538391801782836e-11"
$ws.Range("U16").Value = "[42.90940561697627, 56.02488650670185]"
$ws.Range("Y16").Value = [double]"2.682362362362417"
$ws.Range("Z16").Value = [double]"4.301901901901994"
$ws.Range("M17").Value = "[52.996340850849876, 73.4815583436824]"
$ws.Range("N17").Value = [double]"4.440892098500626e-16"
$ws.Range("O17").Value = [double]"4.440892098500626e-16"
$ws.Range("U17").Value = "[42.18707015145044, 54.090837881448564]"
